# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (tab) right after "总计" and before the
# existing "2021-Q4" sheet, populates it with the Q3-2022 fund holdings,
# and updates the "总计" (summary) sheet with a new row for 2022-Q3 while
# shifting the existing 2021-Q4 / 2020-Q4 rows down.

$wb = $excel.ActiveWorkbook

# --- Locate the existing sheets (before any structural changes) ---------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ421 = $wb.Worksheets.Item("2021-Q4")

# --- Insert the new "2022-Q3" sheet right after "总计" -------------------
$wsQ322 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ322.Name = "2022-Q3"

# Worksheet references taken before Worksheets.Add() can go stale, so
# re-resolve the sheets we still need to touch.
$wsQ421 = $wb.Worksheets.Item("2021-Q4")
$wsQ420 = $wb.Worksheets.Item("2020-Q4")

# Copy the header row + index-column styling from the 2021-Q4 sheet (same
# column layout) so the new sheet matches the workbook's look and feel.
$wsQ421.Range("A1:H2").Copy($wsQ322.Range("A1"))

# --- Populate the new "2022-Q3" sheet with its data ----------------------
$wsQ322.Range("B1").Value = "基金代码"
$wsQ322.Range("C1").Value = "基金名称"
$wsQ322.Range("D1").Value = "基金规模"
$wsQ322.Range("E1").Value = "股票总仓位"
$wsQ322.Range("F1").Value = "仓位占比"
$wsQ322.Range("G1").Value = "持有市值(亿元)"
$wsQ322.Range("H1").Value = "仓位排名"

$wsQ322.Range("A2").Value = 0
$wsQ322.Range("B2").Value = "000049"
$wsQ322.Range("C2").Value = "中银标普全球精选自然资源等权重指数（QDII）A"
$wsQ322.Range("D2").Value = "'0.23"
$wsQ322.Range("E2").Value = "'88.52"
$wsQ322.Range("F2").Value = "'1.01"
$wsQ322.Range("G2").Value = "'0.0023"
$wsQ322.Range("H2").Value = 4

$wsQ322.Range("A3").Copy($wsQ322.Range("A3"))
$wsQ322.Range("A3").Value = 1
$wsQ322.Range("B3").Value = "016334"
$wsQ322.Range("C3").Value = "中银标普全球精选自然资源等权重指数（QDII）C"
$wsQ322.Range("D3").Value = "'0.00"
$wsQ322.Range("E3").Value = "'88.52"
$wsQ322.Range("F3").Value = "'1.01"
$wsQ322.Range("G3").Value = 0
$wsQ322.Range("H3").Value = 4

# Give row 3's index cell (A3) the same style as A2 (it lost its style
# when the row-3 content was written fresh above).
$wsQ322.Range("A2").Copy($wsQ322.Range("A3"))
$wsQ322.Range("A3").Value = 1

# Restore "2020-Q4" as the active/selected tab (adding a sheet made the
# new sheet active).
$wsQ420.Activate()

# --- Update the "总计" (summary) sheet -----------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Existing row 2 (was 2021-Q4) becomes the new 2022-Q3 row.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2

# Existing row 3 (was 2020-Q4) becomes the 2021-Q4 row.
$wsTotal.Range("B3").Value = "2021-Q4"

# New row 4 for 2020-Q4, matching the style of the existing index rows.
$wsTotal.Range("A3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2020-Q4"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0
